$d = $word.ActiveDocument

# 1. Merge "5. Procedure for adopting policies (" + "e.g." + " code style, code check-in
#    steps, documentation" runs (which were split apart by gramStart/gramEnd proofErr
#    markers) into a single run of text, removing the grammar-check markers.
[void]$d.Content.Find.Execute(
    "5. Procedure for adopting policies (e.g. code style, code check-in steps, documentation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5. Procedure for adopting policies (e.g. code style, code check-in steps, documentation",
    2)

# 2. Add a new signature line for Austin VanDenPlas after the last existing
#    "Name: ... Date: ..." paragraph (Cayden Hannon), before the section break.
$end = $d.Content
$end.Collapse(0)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Name: Austin </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>VanDenPlas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Date:3/10</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$end.InsertXML($xml)
